# Update "想去人数" (attendee interest count) figures on the
# "展览" (Exhibition) and "全部类型" (All Types) sheets.
#
# The "全部类型" sheet aggregates rows from every other sheet, so its row
# numbers are offset from "展览" (it has two extra rows mixed in from the
# "演出" sheet). Because several rows share identical old values (e.g. two
# different rows both show 25), the updates are addressed by explicit
# row number per sheet rather than by matching on the old value.

$wb = $excel.ActiveWorkbook

# Row (in column F) -> new value, for the "展览" sheet.
$sheetDisplayUpdates = @{
    2  = 12678
    3  = 611
    5  = 19
    6  = 285
    8  = 234
    9  = 12676
    10 = 26
    11 = 3143
    14 = 7
    15 = 19
    16 = 1197
    18 = 130
    19 = 661
    20 = 2842
    21 = 6122
    23 = 3613
    24 = 218
}

# Row (in column F) -> new value, for the "全部类型" sheet.
$allTypesUpdates = @{
    2  = 12678
    3  = 611
    5  = 19
    6  = 285
    9  = 234
    10 = 12676
    11 = 26
    12 = 3143
    15 = 7
    16 = 19
    17 = 1197
    19 = 130
    20 = 661
    21 = 2842
    23 = 6122
    25 = 3613
    26 = 218
}

$ws = $wb.Worksheets.Item("展览")
foreach ($row in $sheetDisplayUpdates.Keys) {
    $ws.Cells.Item($row, 6).Value = $sheetDisplayUpdates[$row]
}

$ws = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $ws.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
